$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells in column D whose new values look like plain numbers (e.g. "1.00")
# must be forced to Text format first, matching how the source data is
# stored as inline strings (t="inlineStr") rather than numeric cells.
$textCells = @("D5", "D6", "D10", "D11", "D14", "D20", "D23", "D25", "D26", "D27", "D31", "D33", "D35", "D37", "D38", "D39", "D43", "D44", "D47", "D48", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = '68.458.41'
$ws.Range("E2").Value = '  +1.67%  '
$ws.Range("D3").Value = '2.642.58'
$ws.Range("E3").Value = '  +1.48%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '599.53'
$ws.Range("E5").Value = '  +1.20%  '
$ws.Range("D6").Value = '155.04'
$ws.Range("E6").Value = '  +3.15%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +0.49%  '
$ws.Range("D9").Value = '2.642.04'
$ws.Range("E9").Value = '  +1.51%  '
$ws.Range("D10").Value = '0.137'
$ws.Range("E10").Value = '  +6.40%  '
$ws.Range("D11").Value = '0.158'
$ws.Range("E11").Value = '  -0.49%  '
$ws.Range("E12").Value = '  +1.43%  '
$ws.Range("E13").Value = '  +2.18%  '
$ws.Range("D14").Value = '28.11'
$ws.Range("E14").Value = '  +3.09%  '
$ws.Range("E15").Value = '  +3.39%  '
$ws.Range("D16").Value = '3.124.18'
$ws.Range("E16").Value = '  +1.53%  '
$ws.Range("D17").Value = '68.229.29'
$ws.Range("E17").Value = '  +1.58%  '
$ws.Range("D18").Value = '2.643.59'
$ws.Range("E18").Value = '  +1.45%  '
$ws.Range("E19").Value = '  +4.06%  '
$ws.Range("D20").Value = '367.09'
$ws.Range("E20").Value = '  -0.63%  '
$ws.Range("E22").Value = '  +2.76%  '
$ws.Range("D23").Value = '4.89'
$ws.Range("E23").Value = '  +2.63%  '
$ws.Range("E24").Value = '  +5.09%  '
$ws.Range("D25").Value = '73.37'
$ws.Range("E25").Value = '  +0.43%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").Value = '10.02'
$ws.Range("E27").Value = '  +1.11%  '
$ws.Range("E28").Value = '  +6.47%  '
$ws.Range("D29").Value = '2.772.70'
$ws.Range("E29").Value = '  +1.54%  '
$ws.Range("E30").Value = '  -0.41%  '
$ws.Range("D31").Value = '576.58'
$ws.Range("E31").Value = '  -0.16%  '
$ws.Range("E32").Value = '  +5.22%  '
$ws.Range("D33").Value = '8.01'
$ws.Range("E33").Value = '  +4.61%  '
$ws.Range("E34").Value = '  +2.92%  '
$ws.Range("D35").Value = '0.131'
$ws.Range("E35").Value = '  +3.69%  '
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("D37").Value = '1.56'
$ws.Range("E37").Value = '  +3.65%  '
$ws.Range("D38").Value = '160.66'
$ws.Range("E38").Value = '  +2.00%  '
$ws.Range("D39").Value = '19.36'
$ws.Range("E39").Value = '  +1.72%  '
$ws.Range("E40").Value = '  +3.84%  '
$ws.Range("E41").Value = '  +1.22%  '
$ws.Range("E42").Value = '  +3.85%  '
$ws.Range("D43").Value = '2.66'
$ws.Range("E43").Value = '  +4.51%  '
$ws.Range("D44").Value = '17.75'
$ws.Range("E44").Value = '  +3.72%  '
$ws.Range("E45").Value = '  +14.20%  '
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("D47").Value = '40.41'
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("D48").Value = '158.69'
$ws.Range("E49").Value = '  +3.59%  '
$ws.Range("E50").Value = '  +2.88%  '
$ws.Range("D51").Value = '22.06'
$ws.Range("E51").Value = '  +3.54%  '
